$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "items": price/qty tweaks, row 10 restyled into a RAM line, and the
# trailing SSD/RAM/Monitor/Fiver-cable rows (11-13) removed since they moved
# (as sold entries) onto the soldProduct sheet.
# ---------------------------------------------------------------------------
$items = $wb.Worksheets.Item("items")

$items.Range("C3").Value = 100000
$items.Range("D3").Value = 10
$items.Range("D4").Value = 4999

$items.Range("A10").Value = 1016
$items.Range("B10").Value = "RAM"
$items.Range("C10").Value = 3000
$items.Range("D10").Value = 19

$items.Rows("11:13").Delete()

# ---------------------------------------------------------------------------
# Sheet "soldProduct": append five newly recorded sales (rows 25-29).
# ---------------------------------------------------------------------------
$sold = $wb.Worksheets.Item("soldProduct")

$sold.Range("A25").Value = 1016
$sold.Range("B25").Value = "Monitor"
$sold.Range("C25").Value = 2
$sold.Range("D25").Value = 15000
$sold.Range("E25").Value = 30000
$sold.Range("F25").Value = "Partho"
$sold.Range("G25").Value = "fsdjhfg"
$sold.Range("H25").Value = "'54353"
$sold.Range("I25").Value = "COD"

$sold.Range("A26").Value = 1015
$sold.Range("B26").Value = "Ram"
$sold.Range("C26").Value = 2
$sold.Range("D26").Value = 3000
$sold.Range("E26").Value = 6000
$sold.Range("F26").Value = "Partho12"
$sold.Range("G26").Value = "fsdf"
$sold.Range("H26").Value = "'543523"
$sold.Range("I26").Value = "nogod"

$sold.Range("A27").Value = 1014
$sold.Range("B27").Value = "SSD"
$sold.Range("C27").Value = 2
$sold.Range("D27").Value = 2800
$sold.Range("E27").Value = 5600
$sold.Range("F27").Value = "Dipanker"
$sold.Range("G27").Value = "fsdalfk"
$sold.Range("H27").Value = "'43532"
$sold.Range("I27").Value = "COD"

$sold.Range("A28").Value = 1016
$sold.Range("B28").Value = "RAM"
$sold.Range("C28").Value = 1
$sold.Range("D28").Value = 3000
$sold.Range("E28").Value = 3000
$sold.Range("F28").Value = "Safi ahmed"
$sold.Range("G28").Value = "abv"
$sold.Range("H28").Value = "'54325"
$sold.Range("I28").Value = "bkash"

$sold.Range("A29").Value = 1002
$sold.Range("B29").Value = "mouse pad"
$sold.Range("C29").Value = 1
$sold.Range("D29").Value = 200
$sold.Range("E29").Value = 200
$sold.Range("F29").Value = "Safi"
$sold.Range("G29").Value = "fsadf"
$sold.Range("H29").Value = "'32423"
$sold.Range("I29").Value = "COD"

# ---------------------------------------------------------------------------
# Sheet "userAccount": append three newly registered users (rows 10-12).
# ---------------------------------------------------------------------------
$users = $wb.Worksheets.Item("userAccount")

$users.Range("A10").Value = "Saquib"
$users.Range("B10").Value = "saquib"
$users.Range("C10").Value = "'123"
$users.Range("D10").Value = "fsdf"
$users.Range("E10").Value = "'4324"

$users.Range("A11").Value = "Safi"
$users.Range("B11").Value = "safi"
$users.Range("C11").Value = "'123"
$users.Range("D11").Value = "fsdf"
$users.Range("E11").Value = "'44234"

$users.Range("A12").Value = "Saddam"
$users.Range("B12").Value = "saddam"
$users.Range("C12").Value = "'123"
$users.Range("D12").Value = "fsdal;"
$users.Range("E12").Value = "'453"
